$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Range("B3").Value = 2
$wsParams.Range("B12").Value = 4

$wsComp = $wb.Worksheets.Item("comp_quantity_inst1")
$wsComp.Range("A2").Value = "T1"
$wsComp.Range("B2").Value = "T2"
$wsComp.Range("C2").Value = 52
$wsComp.Range("D2").Value = 1

$wsComp.Range("A3").Value = "T2"
$wsComp.Range("B3").Value = "T1"
$wsComp.Range("C3").Value = 29
$wsComp.Range("D3").Value = 1

$wsComp.Range("A4").Value = "T1"
$wsComp.Range("B4").Value = "T2"
$wsComp.Range("C4").Value = 450
$wsComp.Range("D4").Value = 0

$wsComp.Range("A5").Value = "T2"
$wsComp.Range("B5").Value = "T1"
$wsComp.Range("C5").Value = 158
$wsComp.Range("D5").Value = 0

$wsComp.Range("A6:D7").ClearContents()

$wsComp.Range("F7").Select()
